$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column header in H1, copying the formatting of the
# neighboring "sum" header (G1) so it matches the other bold/bordered
# header cells exactly.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the corresponding value for the new column in row 2.
$ws.Range("H2").Value = 1
